# Refresh the cryptocurrency snapshot (Price / Volume(1h) columns)
# with the latest scraped values, mirroring the GitHub Actions job
# that regenerates cryptos.xlsx on a schedule.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Price 58.944.67 -> 58.881.99
$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "58.881.99"
$ws.Range("E2").Value = "  -0.53%  "

# Row 3: Price 2.529.39 -> 2.498.85
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "2.498.85"
$ws.Range("E3").Value = "  +2.23%  "

# Row 4: Price 0.998 -> 1.00
$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "1.00"
$ws.Range("E4").Value = "  +0.12%  "

# Row 5: Price 537.45 -> 536.40
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "536.40"
$ws.Range("E5").Value = "  +0.60%  "

# Row 6: Price 143.37 -> 143.50
$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "143.50"
$ws.Range("E6").Value = "  -2.50%  "

$ws.Range("E7").Value = "  -0.03%  "

# Row 8: Price 0.571 -> 0.570
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "0.570"
$ws.Range("E8").Value = "  +0.59%  "

# Row 9: Price 2.526.67 -> 2.524.43
$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "2.524.43"
$ws.Range("E9").Value = "  +2.76%  "

$ws.Range("E10").Value = "  +0.62%  "

$ws.Range("E11").Value = "  +0.07%  "

# Row 12: Price 5.54 -> 5.53
$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "5.53"
$ws.Range("E12").Value = "  +3.07%  "

$ws.Range("E13").Value = "  -0.18%  "

# Row 14: Price 2.941.36 -> 2.939.94
$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "2.939.94"
$ws.Range("E14").Value = "  +2.19%  "

# Row 15: Price 23.40 -> 23.34
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "23.34"
$ws.Range("E15").Value = "  -3.13%  "

# Row 16: Price 58.868.25 -> 58.800.79
$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "58.800.79"
$ws.Range("E16").Value = "  -0.63%  "

$ws.Range("E17").Value = "  +0.57%  "

# Row 18: Price 2.515.37 -> 2.523.56
$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "2.523.56"
$ws.Range("E18").Value = "  +1.29%  "

# Row 19: Price 11.18 -> 11.17
$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "11.17"
$ws.Range("E19").Value = "  +0.62%  "

$ws.Range("E20").Value = "  -2.49%  "

# Row 21: Price 322.14 -> 322.18
$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "322.18"
$ws.Range("E21").Value = "  -0.56%  "

$ws.Range("E22").Value = "  +3.09%  "

# Row 23: Price 5.74 -> 5.75
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "5.75"
$ws.Range("E23").Value = "  +1.39%  "

# Row 24: Price 61.77 -> 61.71
$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "61.71"
$ws.Range("E24").Value = "  +2.50%  "

# Row 25: Price 0.436 -> 0.435
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "0.435"
$ws.Range("E25").Value = "  -6.80%  "

$ws.Range("E26").Value = "  +0.84%  "

# Row 27: Price 2.614.24 -> 2.624.99
$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "2.624.99"
$ws.Range("E27").Value = "  +2.72%  "

# Row 28: Price 0.993 -> 0.997
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "0.997"
$ws.Range("E28").Value = "  +2.24%  "

$ws.Range("E29").Value = "  +0.22%  "

# Row 30: Price 6.67 -> 6.66
$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "6.66"
$ws.Range("E30").Value = "  -3.58%  "

$ws.Range("E31").Value = "  -1.48%  "

$ws.Range("E32").Value = "  -0.78%  "

$ws.Range("E33").Value = "  -8.67%  "

# Row 34: Price 0.996 -> 0.999
$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "0.999"
$ws.Range("E34").Value = "  +0.09%  "

# Row 35: Price 158.22 -> 157.85
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "157.85"
$ws.Range("E35").Value = "  +0.53%  "

$ws.Range("E36").Value = "  +5.67%  "

# Row 37: Price 18.57 -> 18.56
$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "18.56"
$ws.Range("E37").Value = "  +1.64%  "

$ws.Range("E38").Value = "  -4.50%  "

$ws.Range("E39").Value = "  -6.82%  "

# Row 40: Price 36.40 -> 36.45
$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "36.45"
$ws.Range("E40").Value = "  -0.86%  "

# Row 41: Price 5.52 -> 5.53
$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "5.53"
$ws.Range("E41").Value = "  -4.52%  "

# Row 42: Price 296.58 -> 295.83
$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "295.83"
$ws.Range("E42").Value = "  -5.17%  "

# Row 43: Price 3.64 -> 3.63
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "3.63"
$ws.Range("E43").Value = "  -2.28%  "

$ws.Range("E44").Value = "  -5.35%  "

# Row 45: Price 0.997 -> 0.995
$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "0.995"
$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("E46").Value = "  +3.59%  "

# Row 47: Price 10.78 -> 10.76
$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "10.76"
$ws.Range("E47").Value = "  +0.63%  "

# Row 48: Price 124.59 -> 124.78
$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "124.78"
$ws.Range("E48").Value = "  +4.70%  "

# Row 49: Price 0.0929 -> 0.0928
$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "0.0928"
$ws.Range("E49").Value = "  -0.93%  "

$ws.Range("E50").Value = "  +0.71%  "

# Row 51: Price 0.0512 -> 0.0511
$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"  # keep this a text value, not a number
$dCell.Value = "0.0511"
$ws.Range("E51").Value = "  -1.99%  "
